$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values ---
$ws.Range("B2").Value = "22.02.2024"
$ws.Range("H2").Value = "Document 1800000442 was posted in company code 3000"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1800000442"
$ws.Range("I2").Style = "Normal"

# --- Fill in row 3 (clearing the pre-existing style on A3) ---
$ws.Range("A3").Value = "CMS0000043"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "22.02.2024"
$ws.Range("C3").Value = 1500
$ws.Range("D3").Value = "Robo Test1"
$ws.Range("E3").Value = 40001
$ws.Range("F3").Value = 1500
$ws.Range("G3").Value = 329
$ws.Range("H3").Value = "Document 1800000443 was posted in company code 3000"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1800000443"
$ws.Range("I3").Style = "Normal"

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "3000"
$ws.Range("J3").Style = "Normal"

$ws.Range("K3").Value = 2024

# --- Fill in row 4 (clearing the pre-existing style on A4) ---
$ws.Range("A4").Value = "CMS0000043"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "22.02.2024"
$ws.Range("C4").Value = 2500
$ws.Range("D4").Value = "Robo Test2"
$ws.Range("E4").Value = 40001
$ws.Range("F4").Value = 2500
$ws.Range("G4").Value = 329
$ws.Range("H4").Value = "Document 1800000444 was posted in company code 3000"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1800000444"
$ws.Range("I4").Style = "Normal"

$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "3000"
$ws.Range("J4").Style = "Normal"

$ws.Range("K4").Value = 2024

# --- Update selection to F3 ---
$ws.Range("F3").Select()
